$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, new value) updates, derived from the
# scheduled market-data refresh (currentAveragePrice* / LevePrice* / LeveProfit* columns).

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1599.6666
$ws.Range("I38").Value = 532.6667
$ws.Range("J38").Value = 2666.6667
$ws.Range("K38").Value = 1598.0001
$ws.Range("L38").Value = 8000.000100000001
$ws.Range("M38").Value = -1226.0001
$ws.Range("N38").Value = -8744.000100000001

$ws.Range("H113").Value = 1940.9166
$ws.Range("I113").Value = 1537.35
$ws.Range("J113").Value = 2445.375
$ws.Range("K113").Value = 1537.35
$ws.Range("L113").Value = 2445.375
$ws.Range("M113").Value = 1716.65
$ws.Range("N113").Value = -8953.375

$ws.Range("H129").Value = 1050.3846
$ws.Range("J129").Value = 1139.0857
$ws.Range("L129").Value = 3417.2571
$ws.Range("N129").Value = -13417.2571

$ws.Range("H138").Value = 4834.137
$ws.Range("I138").Value = 543.5185
$ws.Range("J138").Value = 7352.5435
$ws.Range("K138").Value = 1630.5555
$ws.Range("L138").Value = 22057.6305
$ws.Range("M138").Value = 3509.4445
$ws.Range("N138").Value = -32337.6305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 303.5625
$ws.Range("I5").Value = 294.42856
$ws.Range("J5").Value = 367.5
$ws.Range("K5").Value = 294.42856
$ws.Range("L5").Value = 367.5
$ws.Range("M5").Value = -182.42856
$ws.Range("N5").Value = -591.5

$ws.Range("H45").Value = 1784.174
$ws.Range("I45").Value = 1748.5264
$ws.Range("J45").Value = 1953.5
$ws.Range("K45").Value = 1748.5264
$ws.Range("L45").Value = 1953.5
$ws.Range("M45").Value = -1371.5264
$ws.Range("N45").Value = -2707.5

$ws.Range("H102").Value = 2144.3044
$ws.Range("I102").Value = 1854.6666
$ws.Range("J102").Value = 2687.375
$ws.Range("K102").Value = 1854.6666
$ws.Range("L102").Value = 2687.375
$ws.Range("M102").Value = -232.6666
$ws.Range("N102").Value = -5931.375

$ws.Range("H126").Value = 11750
$ws.Range("I126").Value = 11750
$ws.Range("K126").Value = 35250
$ws.Range("M126").Value = -32780

$ws.Range("H132").Value = 1988.7028
$ws.Range("I132").Value = 1450.6086
$ws.Range("J132").Value = 2872.7144
$ws.Range("K132").Value = 4351.825800000001
$ws.Range("L132").Value = 8618.143199999999
$ws.Range("M132").Value = -1821.825800000001
$ws.Range("N132").Value = -13678.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 303.5625
$ws.Range("I4").Value = 294.42856
$ws.Range("J4").Value = 367.5
$ws.Range("K4").Value = 294.42856
$ws.Range("L4").Value = 367.5
$ws.Range("M4").Value = -179.42856
$ws.Range("N4").Value = -597.5

$ws.Range("H20").Value = 1208.6154
$ws.Range("I20").Value = 1142.6666
$ws.Range("K20").Value = 1142.6666
$ws.Range("M20").Value = -895.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4484.0933
$ws.Range("I31").Value = 6098.6665
$ws.Range("K31").Value = 6098.6665
$ws.Range("M31").Value = -5803.6665

$ws.Range("H34").Value = 4484.0933
$ws.Range("I34").Value = 6098.6665
$ws.Range("K34").Value = 6098.6665
$ws.Range("M34").Value = -5896.6665

$ws.Range("H58").Value = 2600142
$ws.Range("I58").Value = 5683921
$ws.Range("J58").Value = 3275.4211
$ws.Range("K58").Value = 5683921
$ws.Range("L58").Value = 3275.4211
$ws.Range("M58").Value = -5683718
$ws.Range("N58").Value = -3681.4211

$ws.Range("H134").Value = 2902.3125
$ws.Range("I134").Value = 1674.5555
$ws.Range("J134").Value = 4480.857
$ws.Range("K134").Value = 5023.666499999999
$ws.Range("L134").Value = 13442.571
$ws.Range("M134").Value = -2488.666499999999
$ws.Range("N134").Value = -18512.571

$ws.Range("H136").Value = 2600142
$ws.Range("I136").Value = 5683921
$ws.Range("J136").Value = 3275.4211
$ws.Range("K136").Value = 17051763
$ws.Range("L136").Value = 9826.263300000001
$ws.Range("M136").Value = -17049213
$ws.Range("N136").Value = -14926.2633

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9259797
$ws.Range("J5").Value = 41667240
$ws.Range("L5").Value = 125001720
$ws.Range("N5").Value = -125001944

$ws.Range("H127").Value = 3233.3125
$ws.Range("J127").Value = 3233.3125
$ws.Range("L127").Value = 9699.9375
$ws.Range("N127").Value = -19619.9375

$ws.Range("H131").Value = 38213.81
$ws.Range("I131").Value = 1732.9166
$ws.Range("J131").Value = 69483.14
$ws.Range("K131").Value = 5198.7498
$ws.Range("L131").Value = 208449.42
$ws.Range("M131").Value = -158.7497999999996
$ws.Range("N131").Value = -218529.42

$ws.Range("H135").Value = 9259797
$ws.Range("J135").Value = 41667240
$ws.Range("L135").Value = 375005160
$ws.Range("N135").Value = -375010230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 24459
$ws.Range("I46").Value = 5415.5
$ws.Range("K46").Value = 5415.5
$ws.Range("M46").Value = -5259.5

$ws.Range("H80").Value = 6218.968
$ws.Range("I80").Value = 8334.846
$ws.Range("J80").Value = 4690.8335
$ws.Range("K80").Value = 8334.846
$ws.Range("L80").Value = 4690.8335
$ws.Range("M80").Value = -7336.846
$ws.Range("N80").Value = -6686.8335

$ws.Range("H83").Value = 6218.968
$ws.Range("I83").Value = 8334.846
$ws.Range("J83").Value = 4690.8335
$ws.Range("K83").Value = 41674.23
$ws.Range("L83").Value = 23454.1675
$ws.Range("M83").Value = -36682.23
$ws.Range("N83").Value = -33438.1675

$ws.Range("H122").Value = 8853.75
$ws.Range("I122").Value = 25503
$ws.Range("J122").Value = 3304
$ws.Range("K122").Value = 76509
$ws.Range("L122").Value = 9912
$ws.Range("M122").Value = -74059
$ws.Range("N122").Value = -14812

$ws.Range("H132").Value = 2896.6316
$ws.Range("I132").Value = 2689.875
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 8069.625
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -5539.625
$ws.Range("N132").Value = -17057.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1153.9474
$ws.Range("I16").Value = 776.0833
$ws.Range("J16").Value = 1801.7142
$ws.Range("K16").Value = 776.0833
$ws.Range("L16").Value = 1801.7142
$ws.Range("M16").Value = -606.0833
$ws.Range("N16").Value = -2141.7142

$ws.Range("H40").Value = 3925.7334
$ws.Range("I40").Value = 3638.6
$ws.Range("J40").Value = 4500
$ws.Range("K40").Value = 3638.6
$ws.Range("L40").Value = 4500
$ws.Range("M40").Value = -3502.6
$ws.Range("N40").Value = -4772

$ws.Range("H42").Value = 29500
$ws.Range("J42").Value = 29500
$ws.Range("L42").Value = 29500
$ws.Range("N42").Value = -30626

$ws.Range("H46").Value = 1199.6666
$ws.Range("I46").Value = 1199
$ws.Range("K46").Value = 1199
$ws.Range("M46").Value = -1011

$ws.Range("H49").Value = 29500
$ws.Range("J49").Value = 29500
$ws.Range("L49").Value = 29500
$ws.Range("N49").Value = -29794

$ws.Range("H136").Value = 6036.8335
$ws.Range("I136").Value = 3391.875
$ws.Range("J136").Value = 6998.636
$ws.Range("K136").Value = 10175.625
$ws.Range("L136").Value = 20995.908
$ws.Range("M136").Value = -7625.625
$ws.Range("N136").Value = -26095.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10176.281
$ws.Range("I136").Value = 9772.357
$ws.Range("J136").Value = 10490.444
$ws.Range("K136").Value = 29317.071
$ws.Range("L136").Value = 31471.332
$ws.Range("M136").Value = -26767.071
$ws.Range("N136").Value = -36571.33199999999

Write-Output "Applied scheduled market-data refresh to 32 leve rows across 8 sheets."
